$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.040.34'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.952.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '378.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('E7').Value = '  -1.13%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.46'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0839'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.416.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.38'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.947.68'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  +5.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.111.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E19').Value = '  -6.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.85%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('E25').Value = '  +3.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.56%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.166'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.111'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.63'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('E35').Value = '  -3.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '33.49'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.55%  '
$ws.Range('E43').Value = '  -3.20%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.14%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.273'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.003.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0332'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.44%  '
